# Remove the empty Title/Content placeholder shapes from the image-only
# slides (they were left over from the layout but never used). Affects
# the slides whose only real content is a full-bleed picture:
#   slide 4 -> all_ppl.png
#   slide 5 -> all_speaker.png
#   slide 7 -> all_questions.png
#   slide 8 -> one_question.png

$p = $ppt.ActivePresentation

$targetSlides = @(4, 5, 7, 8)

foreach ($idx in $targetSlides) {
    $s = $p.Slides.Item($idx)

    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $sh = $s.Shapes.Item($i)
        if ($sh.Name -eq "Title 1" -or $sh.Name -eq "Content Placeholder 2") {
            $sh.Delete()
        }
    }
}
